# Atualização de bases das ligas, do dia: 19-02-2024 às 20:58
#
# The sheet holds one match per row (columns A..AC). This update corrects a
# handful of rows whose HomeTeam/AwayTeam (and all of their associated odds
# columns) had been written to the wrong row - two matches per swapped pair
# simply need their B and F..AC column values exchanged. One further row
# (299) gets its data corrected/completed in place, including three
# previously-missing cells (H/I/J) and two previously-missing PL columns
# (AB/AC).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns B, F..AC (everything describing the fixture except the running
# id in A and the constant Div/Div-Original-Name/Date columns C/D/E).
$swapCols = @("B","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC")

function Swap-Rows($row1, $row2) {
    foreach ($col in $swapCols) {
        $addr1 = "$col$row1"
        $addr2 = "$col$row2"
        $v1 = $ws.Range($addr1).Value2
        $v2 = $ws.Range($addr2).Value2
        $ws.Range($addr1).Value = $v2
        $ws.Range($addr2).Value = $v1
    }
}

# Rows whose HomeTeam/AwayTeam + odds data were transposed between the pair.
Swap-Rows 110 111
Swap-Rows 112 113
Swap-Rows 254 255
Swap-Rows 264 265

# Row 299: corrected match id/date/teams plus the full score/odds line
# (three cells - H/I/J - and the final two PL columns - AB/AC - were
# missing entirely before and are now populated).
$ws.Range("B299").Value = 7824500
$ws.Range("E299").Value = 45340.86458333334
$ws.Range("F299").Value = "Nacional Asuncion"
$ws.Range("G299").Value = "Sportivo Trinidense"
$ws.Range("H299").Value = 2
$ws.Range("I299").Value = 1
$ws.Range("J299").Value = "H"
$ws.Range("K299").Value = 1.909
$ws.Range("L299").Value = 3.25
$ws.Range("M299").Value = 3.75
$ws.Range("N299").Value = 2
$ws.Range("O299").Value = 3.25
$ws.Range("P299").Value = 3.4
$ws.Range("Q299").Value = -0.25
$ws.Range("R299").Value = 1.775
$ws.Range("S299").Value = 2.025
$ws.Range("T299").Value = 2.25
$ws.Range("U299").Value = 1.9
$ws.Range("V299").Value = 1.9
$ws.Range("W299").Value = 1
$ws.Range("X299").Value = -1
$ws.Range("Y299").Value = -1
$ws.Range("Z299").Value = 0.7749999999999999
$ws.Range("AA299").Value = -1
$ws.Range("AB299").Value = 0.8999999999999999
$ws.Range("AC299").Value = -1
